# Insert a new row of "Logan, UT (community)" hydrology data into Sheet2,
# above the existing row 20, shifting the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, pushing existing rows 20-23 down to 21-24.
$ws.Rows(20).Insert()

$ws.Range("A20").Value = "Logan, UT (community)"
$ws.Range("B20").Value = "Logan Hydrology Data"
$ws.Range("C20").Value = "logan_hydrology_final.xlsx"

# Update the selection to match the authored state.
$ws.Range("D20").Select()
